$d = $word.ActiveDocument

# --- Change 1: merge hyperlink text runs (https://gith + u + b.com/... -> https://github.com/...) ---
$old1 = "https://github.com/IsMo167/Cl-0112/tree/main/Proyecto%20Programado%201/Proyecto%20Programado%201"
$new1 = "https://github.com/IsMo167/Cl-0112/tree/main/Proyecto%20Programado%201/Proyecto%20Programado%201"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null
# Restore the hyperlink character style that Find/Replace drops on the merged run.
$hyperlink = $d.Hyperlinks.Item(1)
$hyperlink.Range.Style = "Hipervnculo"

# --- Change 2: "Decisiones de Diseño" paragraph tweak ---
$old2 = "se planteó que el Gato fuese utilizando la posición del arreglo + 1, mientras que el Cuatro en Línea no se sumara uno, además de que Gato"
$new2 = "se planteó que en ambos juegos fuese utilizando la posición del arreglo + 1, además de que Gato"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Change 3: "Decisiones de Implementación" paragraph expansion ---
$old3 = "se utilizaron mucho los ciclos For y While, principalmente for para crear condiciones para saber si se ganó, verificando diagonales, filas y columnas. En la clase Juego Controlador, se utilizó mucho if y else para la utilización de los métodos, llamándolos y usarlos."
$new3 = "Se utilizaron mucho los ciclos For y While, principalmente for para crear condiciones para saber si se ganó, verificando diagonales, filas y columnas. Implementamos métodos que verifican el estado del juego, como ganador() en la clase Gato y esJuegoTerminado() en CuatroEnLinea. Se hizo hincapié en la reutilización de código mediante la creación de métodos que podrían ser llamados múltiples veces en diferentes contextos, como el método mostrarTablero(). En la clase Juego Controlador, se utilizó mucho if y else para la utilización de los métodos, además de un Switch para el menú. Nos aseguramos de documentar cada método y clase con comentarios. Adoptamos un enfoque iterativo, probando cada parte del código a medida que se desarrollaba, lo que nos permitió identificar y corregir errores."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# --- Change 4: "Puntos de Mejora" paragraph rewrite ---
$old4 = "se podría mejorar con el uso de excepciones para que el usuario si o si utilice los formatos y tipos de atributos que se deben utilizar, también, el añadido de una interfaz gráfica para mejor interacción del usuario con el programa, la utilización de más documentación interna por medio de Doxygen, además, de más modos de juego."
$new4 = "Se podría mejorar con el uso de excepciones para que el usuario si o si utilice los formatos y tipos de atributos que se deben utilizar, también, el añadido de una interfaz gráfica para mejor interacción del usuario con el programa. Revisión y optimización el rendimiento del código, especialmente en la lógica de verificación de ganadores y movimientos. Mejor gestión de tiempo para haber tenido el trabajo completo antes, además de más modos de juego."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null
